# scraper.py now fully functional
# - "season" sheet: update the year values in A2:A3, and leave the
#   selection parked on D13.
# - Make "season" the active sheet/tab (it was "player" before).

$wb = $excel.ActiveWorkbook

$seasonWs = $wb.Worksheets.Item("season")
$seasonWs.Range("A2").Value = 1966
$seasonWs.Range("A3").Value = 1967

# Switch the active sheet to "season" and leave the cursor on D13.
$seasonWs.Activate()
$seasonWs.Range("D13").Select()
